$d = $word.ActiveDocument

# 1. "Title:  Prototyping Labs Supervisor" -> "Title:  Prototyping Lab Supervisor"
$d.Content.Find.Execute("Prototyping Labs Supervisor", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Prototyping Lab Supervisor", 2)

# 2. Merge "(if applicable)" runs - normalize text to single contiguous run text.
$d.Content.Find.Execute("(if applicable)", $false, $false, $false, $false, $false,
                         $true, 1, $false, "(if applicable)", 2)

# 3. Merge "FDM 3D Printers (Ultimaker 3 Extended and Stratasys F170)" runs
$d.Content.Find.Execute("FDM 3D Printers (Ultimaker 3 Extended and Stratasys F170)", $false, $false, $false, $false, $false,
                         $true, 1, $false, "FDM 3D Printers (Ultimaker 3 Extended and Stratasys F170)", 2)

# 4. Merge "nano-particles" sentence runs
$d.Content.Find.Execute("Be aware that printing and sanding of finished prints can release nano-particles that are hazardous to health.", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Be aware that printing and sanding of finished prints can release nano-particles that are hazardous to health.", 2)

# 5. Header1: "Prototyping Labs at GIX" -> "Prototyping Lab at GIX"
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)
    if ($hdr.Range.Text -match "Prototyping Labs at GIX") {
        $hdr.Range.Find.Execute("Prototyping Labs at GIX", $false, $false, $false, $false, $false,
                                 $true, 1, $false, "Prototyping Lab at GIX", 2)
    }
}
